$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

$ws.Range("A33").Value = "idrestaurante 1:1 INT AUT_INCRE"
$ws.Range("B33").Value = "idcontato 1:1 INT AUT_INC"
$ws.Range("C33").Value = "idendereco 1:1 INT AUT_INC"
$ws.Range("D33").Value = "idfoto 1:1 INT AUT_INC"
$ws.Range("E33").Value = "idlogin 1:1  INT AUT_INC"
$ws.Range("F33").Value = "idfeedback 1:1 INT AUT_INC"
$ws.Range("A34").Value = "nomerestaurante VARCHAR(100)"
$ws.Range("B34").Value = "telefoneresidencial VARCHAR(15)"
$ws.Range("C34").Value = "logradouro VARCHAR(50)"
$ws.Range("D34").Value = "fotocapa TEXT"
$ws.Range("E34").Value = "usuario VARCHAR(50)"
$ws.Range("F34").Value = "idrestaurante 1:n INT"
$ws.Range("A35").Value = "categoria ENUM"
$ws.Range("B35").Value = "email VARCHA(100)"
$ws.Range("C35").Value = "numero VARCHAR(10)"
$ws.Range("D35").Value = "foto1 TEXT"
$ws.Range("E35").Value = "email VARCHAR(100)"
$ws.Range("F35").Value = "nome VARCHAR(50)"
$ws.Range("A36").Value = "idcontato 1:1 INT"
$ws.Range("B36").Value = "telefonecelular VARCHAR(15)"
$ws.Range("C36").Value = "complemento VARCHAR(50)"
$ws.Range("D36").Value = "foto2 TEXT"
$ws.Range("E36").Value = "senhaVARCHAR(100)"
$ws.Range("F36").Value = "opiniao VARCHAR(100)"
$ws.Range("A37").Value = "idendereco 1:1INT"
$ws.Range("B37").Value = "site VARCHAR(100)"
$ws.Range("C37").Value = "bairro VARCHAR(30)"
$ws.Range("E37").Value = "datacadastro DATE"
$ws.Range("F37").Value = "datacadastro DATE"
$ws.Range("A38").Value = "cnpj VARCHAR(20)"
$ws.Range("C38").Value = "cep VARCHAR(10)"
$ws.Range("F38").Value = "nota INT"
$ws.Range("A39").Value = "descricao TEXT"
$ws.Range("C39").Value = "estado VARCHAR(30)"
$ws.Range("A40").Value = "faixadepreco INT"
$ws.Range("C40").Value = "cidade VARCHAR(20)"
$ws.Range("A41").Value = "idfeedback 1:n INT"
$ws.Range("A42").Value = "datacriacao DAT E"
$ws.Range("A43").Value = "idfoto 1:n INT"
$ws.Range("A44").Value = "horariofuncionamento DATETIME"
$ws.Range("A45").Value = "status VARCHAR(50)"
